$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows used as format sources (already-styled rows from the existing table):
#  - row 891: A=date style, B..F/H=text style, G=empty style (s=2)
#  - row 890: G has text (s=1) -> used to restyle G when a localisation value is present
$formatRowWithEmptyG = 891
$formatRowWithTextG = 890

$newRows = @(
    @(46077, "Yoan Zouma", 70, 8, 10, 5, "Cheville", 4),
    @(46077, "Omar Benyounes", 70, 6, 5, 0, "", 7),
    @(46077, "Yoann Martelat", 70, 7, 6, 3, "Genou", 8),
    @(46077, "Kamal Bafounta", 70, 7, 6, 1, "Cheville genou", 6),
    @(46077, "Maé Clavel", 70, 7, 7, 3, "Tibia", 5),
    @(46077, "Naim Ighbane", 70, 8, 7, 2, "Genou", 6),
    @(46077, "Sofiane Belle", 70, 6, 4, 1, "Ventre", 7),
    @(46077, "Mehdi Boussaid", 70, 6, 7, 0, "", 7),
    @(46077, "Levy Ndoutoume", 70, 8, 8, 0, "", 8),
    @(46077, "Emmanuel Valey", 70, 7, 5, 5, "Cheville", 4),
    @(46077, "Karahali Souaré", 70, 5, 6, 6, "Cheville", 9),
    @(46077, "Theo Owono", 70, 8, 8, 0, "", 10),
    @(46077, "Ilan Ihaddadene", 70, 9, 7, 0, "", 5),
    @(46077, "Naim Dhib", 70, 7, 4, 2, "Psoas", 5),
    @(46077, "Romain Thunet", 70, 7, 5, 0, "", 9),
    @(46078, "Kamal Bafounta", 70, 6, 5, 0, "", 7),
    @(46078, "Naim Ighbane", 70, 6, 7, 2, "Genou", 8),
    @(46078, "Maé Clavel", 70, 7, 8, 6, "Ischio ", 7),
    @(46078, "Omar Benyounes", 70, 4, 5, 0, "", 3),
    @(46078, "Mehdi Boussaid", 70, 7, 7, 0, "", 9),
    @(46078, "Jeremie Laurent", 70, 6, 6, 0, "", 4),
    @(46078, "Yoann Martelat", 70, 6, 6, 3, "Genou", 7),
    @(46078, "Ilan Ihaddadene", 70, 7, 9, 1, "Ampoule", 8),
    @(46078, "Karahali Souaré", 70, 6, 6, 6, "Cheville", 2),
    @(46078, "Theo Owono", 70, 4, 7, 0, "", 3),
    @(46078, "Mattheo Haon", 70, 5, 2, 0, "", 8),
    @(46078, "Romain Thunet", 70, 8, 7, 3, "Coups", 0),
    @(46078, "Nathanael Beta", 70, 5, 5, 0, "", 4)
)

$startRow = 892
$r = $startRow
foreach ($row in $newRows) {
    $date = $row[0]
    $name = $row[1]
    $volume = $row[2]
    $intensite = $row[3]
    $charge = $row[4]
    $fatigue = $row[5]
    $localisation = $row[6]
    $plaisir = $row[7]

    # Copy the full-row number formatting from an existing row so the new
    # row's styles match the table (date format in A, text font in B-H,
    # the "empty localisation" style in G) without minting new styles.
    $ws.Range("A$formatRowWithEmptyG`:I$formatRowWithEmptyG").Copy()
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $volume
    $ws.Cells.Item($r, 4).Value = $intensite
    $ws.Cells.Item($r, 5).Value = $charge
    $ws.Cells.Item($r, 6).Value = $fatigue
    if ($localisation -ne "") {
        # This row has a non-empty "Localisation douleur" - reapply the
        # text-cell format to G (row 890's G cell already carries it).
        $ws.Range("G$formatRowWithTextG").Copy()
        $ws.Range("G$r").PasteSpecial(-4122)
        $ws.Cells.Item($r, 7).Value = $localisation
    }
    $ws.Cells.Item($r, 8).Value = $plaisir
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"

    $r = $r + 1
}

$ws.Range("K921").Select()

"done"
